# Auto-generated edit script: updates LeveProfit/price cells across multiple job sheets
# per the Ultros_Profits.xlsx market-data refresh (scheduled runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1049.3077
$ws.Range("J17").Value = 996.84
$ws.Range("L17").Value = 2990.52
$ws.Range("N17").Value = -3326.52

$ws.Range("H40").Value = 5143.478
$ws.Range("J40").Value = 5876.923
$ws.Range("L40").Value = 5876.923
$ws.Range("N40").Value = -6226.923

$ws.Range("H64").Value = 15211.111
$ws.Range("I64").Value = 9900
$ws.Range("J64").Value = 15875
$ws.Range("K64").Value = 9900
$ws.Range("L64").Value = 15875
$ws.Range("M64").Value = -9652
$ws.Range("N64").Value = -16371

$ws.Range("H67").Value = 15211.111
$ws.Range("I67").Value = 9900
$ws.Range("J67").Value = 15875
$ws.Range("K67").Value = 9900
$ws.Range("L67").Value = 15875
$ws.Range("M67").Value = -9042
$ws.Range("N67").Value = -17591

$ws.Range("H100").Value = 5725.684
$ws.Range("I100").Value = 3043.818
$ws.Range("J100").Value = 9413.25
$ws.Range("K100").Value = 3043.818
$ws.Range("L100").Value = 9413.25
$ws.Range("M100").Value = -2502.818
$ws.Range("N100").Value = -10495.25

$ws.Range("H107").Value = 1451
$ws.Range("I107").Value = 1899.5714
$ws.Range("J107").Value = 553.8570999999999
$ws.Range("K107").Value = 1899.5714
$ws.Range("L107").Value = 553.8570999999999
$ws.Range("M107").Value = 20.42859999999996
$ws.Range("N107").Value = -4393.8571

$ws.Range("H132").Value = 18786.441
$ws.Range("I132").Value = 3835.9546
$ws.Range("J132").Value = 46195.668
$ws.Range("K132").Value = 11507.8638
$ws.Range("L132").Value = 138587.004
$ws.Range("M132").Value = -8977.863799999999
$ws.Range("N132").Value = -143647.004

$ws.Range("H137").Value = 2515.923
$ws.Range("I137").Value = 2062.2173
$ws.Range("J137").Value = 5994.3335
$ws.Range("K137").Value = 6186.651899999999
$ws.Range("L137").Value = 17983.0005
$ws.Range("M137").Value = -3636.651899999999
$ws.Range("N137").Value = -23083.0005

$ws.Range("H138").Value = 2486.658
$ws.Range("I138").Value = 1253.3658
$ws.Range("J138").Value = 3931.3713
$ws.Range("K138").Value = 3760.0974
$ws.Range("L138").Value = 11794.1139
$ws.Range("M138").Value = 1379.9026
$ws.Range("N138").Value = -22074.1139

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14707203
$ws.Range("I32").Value = 15626280
$ws.Range("K32").Value = 15626280
$ws.Range("M32").Value = -15625993

$ws.Range("H63").Value = 14833.667
$ws.Range("I63").Value = 4501
$ws.Range("J63").Value = 20000
$ws.Range("K63").Value = 4501
$ws.Range("L63").Value = 20000
$ws.Range("M63").Value = -3815
$ws.Range("N63").Value = -21372

$ws.Range("H66").Value = 14833.667
$ws.Range("I66").Value = 4501
$ws.Range("J66").Value = 20000
$ws.Range("K66").Value = 22505
$ws.Range("L66").Value = 100000
$ws.Range("M66").Value = -19073
$ws.Range("N66").Value = -106864

$ws.Range("H92").Value = 67777.5
$ws.Range("J92").Value = 67777.5
$ws.Range("L92").Value = 67777.5
$ws.Range("N92").Value = -72769.5

$ws.Range("H102").Value = 4792.4375
$ws.Range("I102").Value = 4711.933
$ws.Range("K102").Value = 4711.933
$ws.Range("M102").Value = -3089.933

$ws.Range("H132").Value = 2701.16
$ws.Range("I132").Value = 2066.5217
$ws.Range("K132").Value = 6199.5651
$ws.Range("M132").Value = -3669.5651

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 50000
$ws.Range("J68").Value = 50000
$ws.Range("L68").Value = 50000
$ws.Range("N68").Value = -51622

$ws.Range("H71").Value = 50000
$ws.Range("J71").Value = 50000
$ws.Range("L71").Value = 150000
$ws.Range("N71").Value = -158112

$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()

$ws.Range("H76").Value = 26387.334
$ws.Range("J76").Value = 26387.334
$ws.Range("L76").Value = 26387.334
$ws.Range("N76").Value = -27017.334

$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()

$ws.Range("H79").Value = 26387.334
$ws.Range("J79").Value = 26387.334
$ws.Range("L79").Value = 26387.334
$ws.Range("N79").Value = -28571.334

$ws.Range("H86").Value = 3500.9312
$ws.Range("I86").Value = 2622.2
$ws.Range("J86").Value = 4442.4287
$ws.Range("K86").Value = 2622.2
$ws.Range("L86").Value = 4442.4287
$ws.Range("M86").Value = -1499.2
$ws.Range("N86").Value = -6688.4287

$ws.Range("H89").Value = 3500.9312
$ws.Range("I89").Value = 2622.2
$ws.Range("J89").Value = 4442.4287
$ws.Range("K89").Value = 13111
$ws.Range("L89").Value = 22212.1435
$ws.Range("M89").Value = -7495
$ws.Range("N89").Value = -33444.14350000001

$ws.Range("H92").Value = 25000
$ws.Range("J92").Value = 25000
$ws.Range("L92").Value = 25000
$ws.Range("N92").Value = -29992

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1614.1875
$ws.Range("I58").Value = 1295.7916
$ws.Range("J58").Value = 2569.375
$ws.Range("K58").Value = 1295.7916
$ws.Range("L58").Value = 2569.375
$ws.Range("M58").Value = -1092.7916
$ws.Range("N58").Value = -2975.375

$ws.Range("H86").Value = 23815.477
$ws.Range("I86").Value = 29864.166
$ws.Range("J86").Value = 15750.556
$ws.Range("K86").Value = 29864.166
$ws.Range("L86").Value = 15750.556
$ws.Range("M86").Value = -28741.166
$ws.Range("N86").Value = -17996.556

$ws.Range("H89").Value = 23815.477
$ws.Range("I89").Value = 29864.166
$ws.Range("J89").Value = 15750.556
$ws.Range("K89").Value = 149320.83
$ws.Range("L89").Value = 78752.78
$ws.Range("M89").Value = -143704.83
$ws.Range("N89").Value = -89984.78

$ws.Range("H99").Value = 18702204
$ws.Range("J99").Value = 40002220
$ws.Range("L99").Value = 40002220
$ws.Range("N99").Value = -40005216

$ws.Range("H126").Value = 18702204
$ws.Range("J126").Value = 40002220
$ws.Range("L126").Value = 120006660
$ws.Range("N126").Value = -120011600

$ws.Range("H136").Value = 1614.1875
$ws.Range("I136").Value = 1295.7916
$ws.Range("J136").Value = 2569.375
$ws.Range("K136").Value = 3887.3748
$ws.Range("L136").Value = 7708.125
$ws.Range("M136").Value = -1337.3748
$ws.Range("N136").Value = -12808.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 27448832
$ws.Range("I4").Value = 27448832
$ws.Range("K4").Value = 82346496
$ws.Range("M4").Value = -82346384

$ws.Range("H14").Value = 48.8
$ws.Range("I14").Value = 48.8
$ws.Range("K14").Value = 146.4
$ws.Range("M14").Value = 26.60000000000002

$ws.Range("H122").Value = 1049.4
$ws.Range("J122").Value = 1100
$ws.Range("L122").Value = 9900
$ws.Range("N122").Value = -14800

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 14285814
$ws.Range("I2").Value = 40.4
$ws.Range("K2").Value = 40.4
$ws.Range("M2").Value = 72.59999999999999

$ws.Range("H70").Value = 500
$ws.Range("I70").Value = 500
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 500
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -230
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 500
$ws.Range("I73").Value = 500
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 500
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = 436
$ws.Range("N73").ClearContents()

$ws.Range("H80").Value = 6804.067
$ws.Range("I80").Value = 3737.5
$ws.Range("K80").Value = 3737.5
$ws.Range("M80").Value = -2739.5

$ws.Range("H83").Value = 6804.067
$ws.Range("I83").Value = 3737.5
$ws.Range("K83").Value = 18687.5
$ws.Range("M83").Value = -13695.5

$ws.Range("H97").Value = 11707.333
$ws.Range("I97").Value = 644.3333
$ws.Range("J97").Value = 33833.332
$ws.Range("K97").Value = 644.3333
$ws.Range("L97").Value = 33833.332
$ws.Range("M97").Value = -148.3333
$ws.Range("N97").Value = -34825.332

$ws.Range("H102").Value = 3063.7
$ws.Range("I102").Value = 2172.8823
$ws.Range("J102").Value = 4228.615
$ws.Range("K102").Value = 2172.8823
$ws.Range("L102").Value = 4228.615
$ws.Range("M102").Value = -550.8823000000002
$ws.Range("N102").Value = -7472.615

$ws.Range("H122").Value = 6598.923
$ws.Range("I122").Value = 5178.6
$ws.Range("K122").Value = 15535.8
$ws.Range("M122").Value = -13085.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9621.643
$ws.Range("I40").Value = 10785.728
$ws.Range("K40").Value = 10785.728
$ws.Range("M40").Value = -10649.728

$ws.Range("H132").Value = 2047.6833
$ws.Range("I132").Value = 1815.3469
$ws.Range("K132").Value = 5446.0407
$ws.Range("M132").Value = -2916.0407

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1165.7742
$ws.Range("I132").Value = 1188.3103
$ws.Range("K132").Value = 3564.9309
$ws.Range("M132").Value = -1034.9309
